$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")
$ws.Activate()

# Insert a new row above the old row 64 (the "toggle(<channel>[,<sn>])" PWM row)
# for the new PHIDGET PWM Frequency command, pushing every row below down by one.
$ws.Range("A64").EntireRow.Insert()

$ws.Range("B64").Value = "frequency(<value>[,<sn>])"
$ws.Range("C64").Value = "PHIDGET PWM Frequency: <value> in Hz"

# Restore the sheet's selection/scroll position to match the post-edit view.
$ws.Range("C63").Select()
